$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values stay as text by pre-formatting the
# affected cells as Text before assigning their string values.
$textCells = @(
    "D5",
    "D6",
    "D12",
    "D14",
    "D15",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D27",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D37",
    "D39",
    "D40",
    "D42",
    "D44",
    "D45",
    "D46",
    "D48",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, matching the source diff.
$ws.Range("D2").Value = "70.768.88"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "3.472.44"
$ws.Range("E3").Value = "  +2.45%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "588.02"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "180.24"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "3.459.12"
$ws.Range("E7").Value = "  +2.20%  "
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +6.19%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "49.54"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("D14").Value = "693.80"
$ws.Range("D15").Value = "8.82"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").Value = "4.019.22"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").Value = "70.658.79"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").Value = "3.466.50"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").Value = "17.92"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").Value = "11.54"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").Value = "0.915"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").Value = "5.50"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").Value = "17.21"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("D27").Value = "2.72"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "34.08"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "8.93"
$ws.Range("E30").Value = "  +2.55%  "
$ws.Range("D31").Value = "7.27"
$ws.Range("E31").Value = "  +4.27%  "
$ws.Range("D32").Value = "3.97"
$ws.Range("E32").Value = "  +10.35%  "
$ws.Range("D33").Value = "577.50"
$ws.Range("E33").Value = "  +3.95%  "
$ws.Range("D34").Value = "11.17"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "59.12"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "3.603.71"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("D39").Value = "0.142"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("D40").Value = "35.71"
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("D41").Value = "0.0₃0741"
$ws.Range("E41").Value = "  +6.24%  "
$ws.Range("D42").Value = "3.38"
$ws.Range("E42").Value = "  +3.07%  "
$ws.Range("E43").Value = "  +2.33%  "
$ws.Range("D44").Value = "0.341"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.37"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0428"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("D48").Value = "1.47"
$ws.Range("E48").Value = "  +3.18%  "
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "134.08"
$ws.Range("E51").Value = "  +0.99%  "
